$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: Row, A (date serial), B (hour), C (value), D (label text)
$data = @(
  ,@(2, 45970, 8, 0.011, "09.11.20258")
  ,@(3, 45970, 9, 0.011, "09.11.20259")
  ,@(4, 45970, 10, 0.297, "09.11.202510")
  ,@(5, 45970, 11, 0.531, "09.11.202511")
  ,@(6, 45970, 12, 0.596, "09.11.202512")
  ,@(7, 45970, 13, 0.564, "09.11.202513")
  ,@(8, 45970, 14, 0.647, "09.11.202514")
  ,@(9, 45970, 15, 0.54, "09.11.202515")
  ,@(10, 45970, 16, 0.513, "09.11.202516")
  ,@(11, 45970, 17, 0.214, "09.11.202517")
  ,@(12, 45970, 18, 0.044, "09.11.202518")
  ,@(13, 45970, 19, 0.011, "09.11.202519")
  ,@(14, 45970, 20, 0.011, "09.11.202520")
  ,@(15, 45970, 21, 0.011, "09.11.202521")
  ,@(16, 45970, 22, 0.011, "09.11.202522")
  ,@(17, 45970, 23, 0.011, "09.11.202523")
  ,@(18, 45970, 24, 0.011, "09.11.202524")
  ,@(19, 45971, 1, 0.011, "10.11.20251")
  ,@(20, 45971, 2, 0.011, "10.11.20252")
  ,@(21, 45971, 3, 0.011, "10.11.20253")
  ,@(22, 45971, 4, 0.011, "10.11.20254")
  ,@(23, 45971, 5, 0.011, "10.11.20255")
  ,@(24, 45971, 6, 0.011, "10.11.20256")
  ,@(25, 45971, 7, 0.011, "10.11.20257")
  ,@(26, 45971, 8, 0.011, "10.11.20258")
  ,@(27, 45971, 9, 0.041, "10.11.20259")
  ,@(28, 45971, 10, 0.174, "10.11.202510")
  ,@(29, 45971, 11, 0.37, "10.11.202511")
  ,@(30, 45971, 12, 0.604, "10.11.202512")
  ,@(31, 45971, 13, 0.589, "10.11.202513")
  ,@(32, 45971, 14, 0.624, "10.11.202514")
  ,@(33, 45971, 15, 0.53, "10.11.202515")
  ,@(34, 45971, 16, 0.386, "10.11.202516")
  ,@(35, 45971, 17, 0.159, "10.11.202517")
  ,@(36, 45971, 18, 0.034, "10.11.202518")
  ,@(37, 45971, 19, 0.011, "10.11.202519")
  ,@(38, 45971, 20, 0.011, "10.11.202520")
  ,@(39, 45971, 21, 0.011, "10.11.202521")
  ,@(40, 45971, 22, 0.011, "10.11.202522")
  ,@(41, 45971, 23, 0.011, "10.11.202523")
  ,@(42, 45971, 24, 0.011, "10.11.202524")
  ,@(43, 45972, 1, 0.011, "11.11.20251")
  ,@(44, 45972, 2, 0.011, "11.11.20252")
  ,@(45, 45972, 3, 0.011, "11.11.20253")
  ,@(46, 45972, 4, 0.011, "11.11.20254")
  ,@(47, 45972, 5, 0.011, "11.11.20255")
  ,@(48, 45972, 6, 0.011, "11.11.20256")
  ,@(49, 45972, 7, 0.011, "11.11.20257")
  ,@(50, 45972, 8, 0.011, "11.11.20258")
  ,@(51, 45972, 9, 0.024, "11.11.20259")
  ,@(52, 45972, 10, 0.147, "11.11.202510")
  ,@(53, 45972, 11, 0.28, "11.11.202511")
  ,@(54, 45972, 12, 0.44, "11.11.202512")
  ,@(55, 45972, 13, 0.455, "11.11.202513")
  ,@(56, 45972, 14, 0.446, "11.11.202514")
  ,@(57, 45972, 15, 0.236, "11.11.202515")
  ,@(58, 45972, 16, 0.156, "11.11.202516")
  ,@(59, 45972, 17, 0.082, "11.11.202517")
  ,@(60, 45972, 18, 0.016, "11.11.202518")
  ,@(61, 45972, 19, 0.011, "11.11.202519")
  ,@(62, 45972, 20, 0.011, "11.11.202520")
  ,@(63, 45972, 21, 0.011, "11.11.202521")
  ,@(64, 45972, 22, 0.011, "11.11.202522")
  ,@(65, 45972, 23, 0.011, "11.11.202523")
  ,@(66, 45972, 24, 0.011, "11.11.202524")
  ,@(67, 45973, 1, 0.011, "12.11.20251")
  ,@(68, 45973, 2, 0.011, "12.11.20252")
  ,@(69, 45973, 3, 0.011, "12.11.20253")
  ,@(70, 45973, 4, 0.011, "12.11.20254")
  ,@(71, 45973, 5, 0.011, "12.11.20255")
  ,@(72, 45973, 6, 0.011, "12.11.20256")
  ,@(73, 45973, 7, 0.011, "12.11.20257")
  ,@(74, 45973, 8, 0.011, "12.11.20258")
  ,@(75, 45973, 9, 0.033, "12.11.20259")
  ,@(76, 45973, 10, 0.272, "12.11.202510")
  ,@(77, 45973, 11, 0.493, "12.11.202511")
  ,@(78, 45973, 12, 0.624, "12.11.202512")
  ,@(79, 45973, 13, 0.584, "12.11.202513")
  ,@(80, 45973, 14, 0.61, "12.11.202514")
  ,@(81, 45973, 15, 0.644, "12.11.202515")
  ,@(82, 45973, 16, 0.472, "12.11.202516")
  ,@(83, 45973, 17, 0.223, "12.11.202517")
  ,@(84, 45973, 18, 0.034, "12.11.202518")
  ,@(85, 45973, 19, 0.011, "12.11.202519")
  ,@(86, 45973, 20, 0.011, "12.11.202520")
  ,@(87, 45973, 21, 0.011, "12.11.202521")
  ,@(88, 45973, 22, 0.011, "12.11.202522")
  ,@(89, 45973, 23, 0.011, "12.11.202523")
  ,@(90, 45973, 24, 0.011, "12.11.202524")
  ,@(91, 45974, 1, 0.011, "13.11.20251")
  ,@(92, 45974, 2, 0.011, "13.11.20252")
  ,@(93, 45974, 3, 0.011, "13.11.20253")
  ,@(94, 45974, 4, 0.011, "13.11.20254")
  ,@(95, 45974, 5, 0.011, "13.11.20255")
  ,@(96, 45974, 6, 0.011, "13.11.20256")
  ,@(97, 45974, 7, 0.011, "13.11.20257")
  ,@(98, 45974, 8, 0.011, "13.11.20258")
  ,@(99, 45974, 9, 0.13, "13.11.20259")
  ,@(100, 45974, 10, 0.597, "13.11.202510")
  ,@(101, 45974, 11, 0.859, "13.11.202511")
  ,@(102, 45974, 12, 2.172, "13.11.202512")
  ,@(103, 45974, 13, 2.513, "13.11.202513")
  ,@(104, 45974, 14, 2.712, "13.11.202514")
  ,@(105, 45974, 15, 2.489, "13.11.202515")
  ,@(106, 45974, 16, 1.464, "13.11.202516")
  ,@(107, 45974, 17, 0.448, "13.11.202517")
  ,@(108, 45974, 18, 0.067, "13.11.202518")
  ,@(109, 45974, 19, 0.011, "13.11.202519")
  ,@(110, 45974, 20, 0.011, "13.11.202520")
  ,@(111, 45974, 21, 0.011, "13.11.202521")
  ,@(112, 45974, 22, 0.011, "13.11.202522")
  ,@(113, 45974, 23, 0.011, "13.11.202523")
  ,@(114, 45974, 24, 0.011, "13.11.202524")
  ,@(115, 45975, 1, 0.011, "14.11.20251")
  ,@(116, 45975, 2, 0.011, "14.11.20252")
  ,@(117, 45975, 3, 0.011, "14.11.20253")
  ,@(118, 45975, 4, 0.011, "14.11.20254")
  ,@(119, 45975, 5, 0.011, "14.11.20255")
  ,@(120, 45975, 6, 0.011, "14.11.20256")
  ,@(121, 45975, 7, 0.011, "14.11.20257")
  ,@(122, 45975, 8, 0.011, "14.11.20258")
  ,@(123, 45975, 9, 0.138, "14.11.20259")
  ,@(124, 45975, 10, 0.681, "14.11.202510")
  ,@(125, 45975, 11, 2.036, "14.11.202511")
  ,@(126, 45975, 12, 2.586, "14.11.202512")
  ,@(127, 45975, 13, 2.733, "14.11.202513")
  ,@(128, 45975, 14, 2.726, "14.11.202514")
  ,@(129, 45975, 15, 2.316, "14.11.202515")
  ,@(130, 45975, 16, 0.971, "14.11.202516")
  ,@(131, 45975, 17, 0.463, "14.11.202517")
  ,@(132, 45975, 18, 0.064, "14.11.202518")
  ,@(133, 45975, 19, 0.011, "14.11.202519")
  ,@(134, 45975, 20, 0.011, "14.11.202520")
  ,@(135, 45975, 21, 0.011, "14.11.202521")
  ,@(136, 45975, 22, 0.011, "14.11.202522")
  ,@(137, 45975, 23, 0.011, "14.11.202523")
  ,@(138, 45975, 24, 0.011, "14.11.202524")
  ,@(139, 45976, 1, 0.011, "15.11.20251")
  ,@(140, 45976, 2, 0.011, "15.11.20252")
  ,@(141, 45976, 3, 0.011, "15.11.20253")
  ,@(142, 45976, 4, 0.011, "15.11.20254")
  ,@(143, 45976, 5, 0.011, "15.11.20255")
  ,@(144, 45976, 6, 0.011, "15.11.20256")
  ,@(145, 45976, 7, 0.011, "15.11.20257")
  ,@(146, 45976, 8, 0.011, "15.11.20258")
  ,@(147, 45976, 9, 0.13, "15.11.20259")
  ,@(148, 45976, 10, 0.927, "15.11.202510")
  ,@(149, 45976, 11, 2.406, "15.11.202511")
  ,@(150, 45976, 12, 2.913, "15.11.202512")
  ,@(151, 45976, 13, 3.181, "15.11.202513")
  ,@(152, 45976, 14, 3.145, "15.11.202514")
  ,@(153, 45976, 15, 2.664, "15.11.202515")
  ,@(154, 45976, 16, 2.115, "15.11.202516")
  ,@(155, 45976, 17, 0.571, "15.11.202517")
  ,@(156, 45976, 18, 0.065, "15.11.202518")
  ,@(157, 45976, 19, 0.011, "15.11.202519")
  ,@(158, 45976, 20, 0.011, "15.11.202520")
  ,@(159, 45976, 21, 0.011, "15.11.202521")
  ,@(160, 45976, 22, 0.011, "15.11.202522")
  ,@(161, 45976, 23, 0.011, "15.11.202523")
  ,@(162, 45976, 24, 0.011, "15.11.202524")
  ,@(163, 45977, 1, 0.011, "16.11.20251")
  ,@(164, 45977, 2, 0.011, "16.11.20252")
  ,@(165, 45977, 3, 0.011, "16.11.20253")
  ,@(166, 45977, 4, 0.011, "16.11.20254")
  ,@(167, 45977, 5, 0.011, "16.11.20255")
  ,@(168, 45977, 6, 0.011, "16.11.20256")
  ,@(169, 45977, 7, 0.011, "16.11.20257")
  ,@(170, 45977, 8, 0.011, "16.11.20258")
)

foreach ($row in $data) {
  $r = $row[0]
  $a = $row[1]
  $b = $row[2]
  $c = $row[3]
  $d = $row[4]
  $ws.Cells.Item($r, 1).Value = $a
  $ws.Cells.Item($r, 2).Value = $b
  $ws.Cells.Item($r, 3).Value = $c
  $ws.Cells.Item($r, 4).Value = $d
}

Write-Host "Done updating" $data.Count "rows"
